# Insert a new row of weekly price data for "Achicoria" (La Araucanía / Vega
# Modelo de Temuco) at row 21, pushing the existing rows 21-48 down to 22-49.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 21 (shifts rows 21..48 down to 22..49)
$ws.Rows.Item(21).Insert()

# Populate the newly inserted row with the new weekly record
$ws.Range("A21").Value = 10
$ws.Range("B21").Value = "Vega Modelo de Temuco"
$ws.Range("C21").Value = "La Araucanía"
$ws.Range("D21").Value = 44757
$ws.Range("E21").Value = 9
$ws.Range("F21").Value = 100112010
$ws.Range("G21").Value = "Achicoria"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 100
$ws.Range("K21").Value = 10000
$ws.Range("L21").Value = 10000
$ws.Range("M21").Value = 10000
$ws.Range("N21").Value = "$/caja 18 unidades"
$ws.Range("O21").Value = "Región Metropolitana"
$ws.Range("P21").Value = 556
$ws.Range("Q21").Value = 18
$ws.Range("R21").Value = "Hortaliza"
